$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 549
$ws.Range("F5").Value = 14
$ws.Range("F6").Value = 708
$ws.Range("F7").Value = 320
$ws.Range("F9").Value = 117
$ws.Range("F10").Value = 232
$ws.Range("F12").Value = 4568
$ws.Range("F13").Value = 34
$ws.Range("F15").Value = 464
$ws.Range("F17").Value = 522
$ws.Range("F18").Value = 313
$ws.Range("F22").Value = 692
$ws.Range("F23").Value = 79
$ws.Range("F24").Value = 286
$ws.Range("F25").Value = 983
$ws.Range("F26").Value = 58
$ws.Range("F27").Value = 1669
$ws.Range("F28").Value = 395
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 36
$ws.Range("F5").Value = 254
$ws.Range("F6").Value = 35
$ws.Range("F8").Value = 290
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 155
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 36
$ws.Range("F6").Value = 155
$ws.Range("F7").Value = 549
$ws.Range("F9").Value = 14
$ws.Range("F10").Value = 708
$ws.Range("F12").Value = 320
$ws.Range("F14").Value = 117
$ws.Range("F15").Value = 232
$ws.Range("F18").Value = 4569
$ws.Range("F19").Value = 34
$ws.Range("F21").Value = 254
$ws.Range("F22").Value = 464
$ws.Range("F24").Value = 522
$ws.Range("F25").Value = 313
$ws.Range("F27").Value = 35
$ws.Range("F31").Value = 290
$ws.Range("F33").Value = 692
$ws.Range("F37").Value = 79
$ws.Range("F38").Value = 286
$ws.Range("F39").Value = 983
$ws.Range("F40").Value = 58
$ws.Range("F41").Value = 1669
$ws.Range("F42").Value = 395
